$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()
$ws.Range("A2").Value = "Ayati Arvind"
$ws.Range("A2").Select()
